$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the obsolete "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Keep the first sheet active/selected, matching the original workbook view
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
